# Update the iServ stats for the latest month (row 27) to reflect 2026-02 data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Raw source values that changed
$ws.Range("B27").Value = 6548        # schools
$ws.Range("D27").Value = 6105171     # users

# Derived values recomputed from the raw values (workbook stores values, not formulas)
$ws.Range("E27").Value = 932.3718692730605   # users_per_school = users / schools
$ws.Range("F27").Value = 10.05042016806723   # yoy_schools
$ws.Range("H27").Value = 25.22125391267      # yoy_users
